# Auto-generated edit script: updates market-price derived columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets per scheduled refresh.
$wb = $excel.ActiveWorkbook

# Flat triples: SheetName, CellRef, NewValue
$setOps = @(
  "ALC","H126",45000,
  "ALC","J126",45000,
  "ALC","L126",45000,
  "ALC","N126",-54880,
  "ALC","H129",1081.3334,
  "ALC","I129",820.2,
  "ALC","J129",1127.9642,
  "ALC","K129",2460.6,
  "ALC","L129",3383.8926,
  "ALC","M129",2539.4,
  "ALC","N129",-13383.8926,
  "ALC","H137",1284.1632,
  "ALC","I137",1066.3429,
  "ALC","K137",3199.0287,
  "ALC","M137",-649.0287000000003,
  "ARM","H32",15800.564,
  "ARM","I32",15800.564,
  "ARM","K32",15800.564,
  "ARM","M32",-15513.564,
  "ARM","H45",1054.8572,
  "ARM","I45",874,
  "ARM","K45",874,
  "ARM","M45",-497,
  "ARM","H61",2829.8462,
  "ARM","I61",2638.8,
  "ARM","K61",2638.8,
  "ARM","M61",-2426.8,
  "ARM","H74",1415,
  "ARM","I74",1237.5,
  "ARM","K74",1237.5,
  "ARM","M74",-363.5,
  "ARM","H77",1415,
  "ARM","I77",1237.5,
  "ARM","K77",6187.5,
  "ARM","M77",-1819.5,
  "ARM","H132",1884.8292,
  "ARM","I132",1483.6129,
  "ARM","J132",3128.6,
  "ARM","K132",4450.8387,
  "ARM","L132",9385.799999999999,
  "ARM","M132",-1920.8387,
  "ARM","N132",-14445.8,
  "ARM","H136",2829.8462,
  "ARM","I136",2638.8,
  "ARM","K136",7916.400000000001,
  "ARM","M136",-5366.400000000001,
  "BSM","H12",29602,
  "BSM","I12",17336.666,
  "BSM","J12",48000,
  "BSM","K12",17336.666,
  "BSM","L12",48000,
  "BSM","M12",-17168.666,
  "BSM","N12",-48336,
  "BSM","H29",2903.2,
  "BSM","J29",10000,
  "BSM","L29",10000,
  "BSM","N29",-10578,
  "BSM","H134",1709.25,
  "BSM","I134",1595.9375,
  "BSM","J134",2162.5,
  "BSM","K134",4787.8125,
  "BSM","L134",6487.5,
  "BSM","M134",-2252.8125,
  "BSM","N134",-11557.5,
  "CRP","H31",1981.3726,
  "CRP","I31",1401.4872,
  "CRP","J31",3866,
  "CRP","K31",1401.4872,
  "CRP","L31",3866,
  "CRP","M31",-1106.4872,
  "CRP","N31",-4456,
  "CRP","H34",1981.3726,
  "CRP","I34",1401.4872,
  "CRP","J34",3866,
  "CRP","K34",1401.4872,
  "CRP","L34",3866,
  "CRP","M34",-1199.4872,
  "CRP","N34",-4270,
  "CRP","H58",2322.6667,
  "CRP","I58",2490.5,
  "CRP","J58",980,
  "CRP","K58",2490.5,
  "CRP","L58",980,
  "CRP","M58",-2287.5,
  "CRP","N58",-1386,
  "CRP","H64",0,
  "CRP","J64",0,
  "CRP","L64",0,
  "CRP","H67",0,
  "CRP","J67",0,
  "CRP","L67",0,
  "CRP","H74",29438,
  "CRP","J74",29438,
  "CRP","L74",29438,
  "CRP","N74",-31186,
  "CRP","H77",29438,
  "CRP","J77",29438,
  "CRP","L77",88314,
  "CRP","N77",-97050,
  "CRP","H107",397.27777,
  "CRP","I107",440.6154,
  "CRP","J107",284.6,
  "CRP","K107",440.6154,
  "CRP","L107",284.6,
  "CRP","M107",1479.3846,
  "CRP","N107",-4124.6,
  "CRP","H136",2322.6667,
  "CRP","I136",2490.5,
  "CRP","J136",980,
  "CRP","K136",7471.5,
  "CRP","L136",2940,
  "CRP","M136",-4921.5,
  "CRP","N136",-8040,
  "CUL","H11",283.72223,
  "CUL","I11",80.7,
  "CUL","J11",537.5,
  "CUL","K11",242.1,
  "CUL","L11",1612.5,
  "CUL","M11",-102.1,
  "CUL","N11",-1892.5,
  "CUL","H131",17242880,
  "CUL","I131",415,
  "CUL","J131",17858684,
  "CUL","K131",1245,
  "CUL","L131",53576052,
  "CUL","M131",3795,
  "CUL","N131",-53586132,
  "GSM","H41",1051,
  "GSM","I41",1051,
  "GSM","K41",1051,
  "GSM","M41",-696,
  "GSM","H92",26593.334,
  "GSM","J92",26593.334,
  "GSM","L92",26593.334,
  "GSM","N92",-30337.334,
  "GSM","H95",1443406.2,
  "GSM","J95",1443406.2,
  "GSM","L95",1443406.2,
  "GSM","N95",-1448898.2,
  "GSM","H107",441.17856,
  "GSM","I107",260.06668,
  "GSM","J107",650.1539,
  "GSM","K107",260.06668,
  "GSM","L107",650.1539,
  "GSM","M107",1659.93332,
  "GSM","N107",-4490.1539,
  "GSM","H110",49851,
  "GSM","J110",49851,
  "GSM","L110",49851,
  "GSM","N110",-58031,
  "GSM","H132",3373.4119,
  "GSM","I132",2846.25,
  "GSM","J132",3842,
  "GSM","K132",8538.75,
  "GSM","L132",11526,
  "GSM","M132",-6008.75,
  "GSM","N132",-16586,
  "LTW","H16",1398.3334,
  "LTW","I16",1178,
  "LTW","J16",2500,
  "LTW","K16",1178,
  "LTW","L16",2500,
  "LTW","M16",-1008,
  "LTW","N16",-2840,
  "LTW","H64",0,
  "LTW","J64",0,
  "LTW","L64",0,
  "LTW","H67",0,
  "LTW","J67",0,
  "LTW","L67",0,
  "LTW","H136",3080.5,
  "LTW","I136",2350.3125,
  "LTW","K136",7050.9375,
  "LTW","M136",-4500.9375,
  "WVR","H63",500000,
  "WVR","J63",500000,
  "WVR","L63",500000,
  "WVR","N63",-501248,
  "WVR","H66",500000,
  "WVR","J66",500000,
  "WVR","L66",1500000,
  "WVR","N66",-1506240,
  "WVR","H107",725.2353000000001,
  "WVR","I107",694.53845,
  "WVR","J107",825,
  "WVR","K107",2083.61535,
  "WVR","L107",2475,
  "WVR","M107",-163.61535,
  "WVR","N107",-6315,
  "WVR","H122",9261786,
  "WVR","I122",11906881,
  "WVR","K122",35720643,
  "WVR","M122",-35718193,
  "WVR","H126",17248.818,
  "WVR","I126",17248.818,
  "WVR","J126",0,
  "WVR","K126",51746.454,
  "WVR","L126",0,
  "WVR","M126",-49276.454,
  "WVR","H136",2123.6875,
  "WVR","I136",2206.1538,
  "WVR","J136",1766.3334,
  "WVR","K136",6618.4614,
  "WVR","L136",5299.0002,
  "WVR","M136",-4068.4614,
  "WVR","N136",-10399.0002
)

for ($i = 0; $i -lt $setOps.Count; $i += 3) {
  $sheetName = $setOps[$i]
  $cellRef   = $setOps[$i + 1]
  $newValue  = $setOps[$i + 2]
  $ws = $wb.Worksheets.Item($sheetName)
  $ws.Range($cellRef).Value = $newValue
}

# Flat pairs: SheetName, CellRef to clear entirely (cell removed, not zeroed)
$clearOps = @(
  "CRP","N64",
  "CRP","N67",
  "LTW","N64",
  "LTW","N67",
  "WVR","N126"
)

for ($i = 0; $i -lt $clearOps.Count; $i += 2) {
  $sheetName = $clearOps[$i]
  $cellRef   = $clearOps[$i + 1]
  $ws = $wb.Worksheets.Item($sheetName)
  $ws.Range($cellRef).ClearContents()
}

Write-Output "Scheduled runner update complete."
